$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/percentage cells (Coin name, Link, Volume) - safe to assign directly.
$textUpdates = @{
    'E2' = '  +0.49%  '
    'E3' = '  -0.44%  '
    'E4' = '  -0.09%  '
    'E5' = '  +0.24%  '
    'E6' = '  -0.07%  '
    'E7' = '  -1.64%  '
    'E8' = '  -0.64%  '
    'E9' = '  -2.01%  '
    'E10' = '  -1.57%  '
    'E11' = '  +0.38%  '
    'E12' = '  +1.12%  '
    'E13' = '  -0.97%  '
    'E14' = '  -1.57%  '
    'E15' = '  -1.94%  '
    'E16' = '  +0.35%  '
    'E17' = '  -0.08%  '
    'E18' = '  +2.25%  '
    'E19' = '  +0.21%  '
    'E20' = '  -3.10%  '
    'E21' = '  -0.03%  '
    'E22' = '  -3.36%  '
    'E23' = '  -2.48%  '
    'B24' = 'Monero'
    'C24' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'E24' = '  -1.08%  '
    'B25' = 'Cosmos'
    'C25' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    'E25' = '  -0.80%  '
    'E26' = '  +0.23%  '
    'E27' = '  +0.63%  '
    'E28' = '  -0.50%  '
    'E29' = '  -0.54%  '
    'E30' = '  -2.92%  '
    'E31' = '  -2.15%  '
    'E32' = '  -0.10%  '
    'E33' = '  -0.35%  '
    'E34' = '  -2.84%  '
    'E35' = '  +0.13%  '
    'E36' = '  +0.79%  '
    'E37' = '  +1.30%  '
    'E38' = '  -4.52%  '
    'E39' = '  +0.37%  '
    'E40' = '  -2.31%  '
    'E42' = '  -0.16%  '
    'E43' = '  -1.75%  '
    'E44' = '  +0.70%  '
    'E45' = '  +0.86%  '
    'E46' = '  -0.11%  '
    'E47' = '  -2.00%  '
    'E48' = '  +1.21%  '
    'E49' = '  -1.12%  '
    'E50' = '  -3.83%  '
    'E51' = '  -1.60%  '
}
foreach ($cell in $textUpdates.Keys) {
    $ws.Range($cell).Value = $textUpdates[$cell]
}

# Price cells (column D) look numeric (e.g. "0.9997", "4.938", "30.417.83").
# Assigning such strings straight to .Value lets Excel auto-convert them to
# real numbers, which would change the stored cell type from text to numeric.
# The source workbook stores these as plain text, so round-trip the new value
# through a Text-formatted helper cell and Paste Special > Values: that copies
# the text typed value across without touching the destination cell style.
$priceUpdates = @{
    'D2' = '30.417.83'
    'D3' = '1.848.93'
    'D4' = '0.9997'
    'D5' = '233.41'
    'D7' = '0.4672'
    'D8' = '0.2733'
    'D9' = '0.06294'
    'D10' = '1.824.09'
    'D11' = '0.07456'
    'D12' = '16.27'
    'D13' = '4.938'
    'D14' = '83.97'
    'D15' = '0.6204'
    'D16' = '30.358.35'
    'D17' = '0.9996'
    'D18' = '229.78'
    'D19' = '0.000007331'
    'D21' = '1.001'
    'D22' = '4.924'
    'D23' = '5.868'
    'D24' = '165.55'
    'D25' = '9.148'
    'D26' = '17.84'
    'D28' = '0.1023'
    'D30' = '4.094'
    'D31' = '3.814'
    'D32' = '0.04877'
    'D34' = '0.7041'
    'D35' = '2.687'
    'D36' = '0.01918'
    'D37' = '2.663'
    'D38' = '0.8620'
    'D39' = '105.76'
    'D40' = '1.930'
    'D42' = '5.520'
    'D43' = '0.4028'
    'D44' = '7.085'
    'D45' = '61.69'
    'D47' = '8.636'
    'D48' = '33.36'
    'D49' = '0.05529'
    'D50' = '1.348'
    'D51' = '0.3649'
}
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"
foreach ($cell in $priceUpdates.Keys) {
    $helper.Value = $priceUpdates[$cell]
    $helper.Copy()
    $ws.Range($cell).PasteSpecial(-4163)
}
$helper.Clear()
$excel.CutCopyMode = 0

